$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-37 down to 6-38
$ws.Rows(5).Insert()

# Populate the newly inserted row 5 with the new data record
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C5").Value = 'Metropolitana'
$ws.Range("D5").Value = [DateTime]"2022-02-21"
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = 'Tropicales y subtropicales'
$ws.Range("I5").Value = 100108007
$ws.Range("J5").Value = 'Coco'
$ws.Range("K5").Value = 'Sin especificar'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("Q5").Value = '$/malla 20 unidades'
$ws.Range("R5").Value = 'Perú'
$ws.Range("S5").Value = 1500
$ws.Range("T5").Value = 20
